# Append two new listings scraped at 2025-11-03 18:25:06 JST to the
# "ランサーズ" sheet, and refresh the "取得日時" timestamp on the rows
# that were already present from the previous run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-03 18:25:06"

# Refresh the capture timestamp on the existing rows (2-4).
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(4, 1).Value = $newTimestamp

# Row 5: new listing
$ws.Cells.Item(5, 1).Value = $newTimestamp
$ws.Cells.Item(5, 2).Value = "【急募】LINEオプチャ指示を基にしたMT4自動発注システム構築"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5426185"
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5426185")
$ws.Cells.Item(5, 6).Style = "Hyperlink"
$ws.Cells.Item(5, 7).Value = 33

# Row 6: new listing
$ws.Cells.Item(6, 1).Value = $newTimestamp
$ws.Cells.Item(6, 2).Value = "【緊急募集】Laravel + Vue.js オンライン日本語学校システムのバグ修正エンジニア"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5426038"
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5426038")
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Cells.Item(6, 7).Value = 28
